# Sync automático del tracker (cada 3h)
# Appends new match rows (93-103) to the tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=93;  EventId=14357969; Fecha="2025-08-07"; JugadorA="Borna Ćorić";              JugadorB="Christopher O'Connell"; Pronostico="Gana Christopher O'Connell";  Cuota=2.38 },
    @{ Row=94;  EventId=14357970; Fecha="2025-08-07"; JugadorA="Benjamin Bonzi";           JugadorB="Matteo Arnaldi";         Pronostico="Gana Benjamin Bonzi";          Cuota=3 },
    @{ Row=95;  EventId=14357977; Fecha="2025-08-07"; JugadorA="Tomas Martin Etcheverry";  JugadorB="Juncheng Shang";         Pronostico="Gana Juncheng Shang";          Cuota=2.3 },
    @{ Row=96;  EventId=14358008; Fecha="2025-08-07"; JugadorA="Corentin Moutet";          JugadorB="Mackenzie McDonald";     Pronostico="Gana Corentin Moutet";         Cuota=1.67 },
    @{ Row=97;  EventId=14357971; Fecha="2025-08-08"; JugadorA="Sebastián Báez";           JugadorB="David Goffin";           Pronostico="Gana Sebastián Báez";          Cuota=2.75 },
    @{ Row=98;  EventId=14366991; Fecha="2025-08-07"; JugadorA="Sorana Cirstea";           JugadorB="Donna Vekić";            Pronostico="Gana Donna Vekić";             Cuota=1.57 },
    @{ Row=99;  EventId=14366989; Fecha="2025-08-08"; JugadorA="Peyton Stearns";           JugadorB="Yafan Wang";             Pronostico="Gana Peyton Stearns";          Cuota=1.53 },
    @{ Row=100; EventId=14311737; Fecha="2025-08-07"; JugadorA="Petr Brunclik";            JugadorB="Giulio Zeppieri";        Pronostico="Gana Petr Brunclik";           Cuota=3.5 },
    @{ Row=101; EventId=14311734; Fecha="2025-08-08"; JugadorA="Francesco Maestrelli";     JugadorB="Dino Prižmić";           Pronostico="Gana Francesco Maestrelli";    Cuota=3.75 },
    @{ Row=102; EventId=14311725; Fecha="2025-08-08"; JugadorA="Kamil Majchrzak";          JugadorB="Daniil Glinka";          Pronostico="Gana Daniil Glinka";           Cuota=4 },
    @{ Row=103; EventId=14310261; Fecha="2025-08-07"; JugadorA="Giles Hussey";             JugadorB="Stefan Dostanic";        Pronostico="Gana Giles Hussey";            Cuota=2.2 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.EventId

    # Force the date-like text to stay text instead of auto-converting to a
    # date serial number, then drop back to the default style so no
    # formatting residue is left on the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $r.Fecha
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $r.JugadorA
    $ws.Cells.Item($row, 4).Value = $r.JugadorB
    $ws.Cells.Item($row, 5).Value = $r.Pronostico
    $ws.Cells.Item($row, 6).Value = $r.Cuota
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = ""
}
